$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 54
$ws.Range("H54").Value = 24499.666
$ws.Range("I54").Value = 24250
$ws.Range("K54").Value = 24250
$ws.Range("M54").Value = -23764

# Row 55
$ws.Range("H55").Value = 924.6
$ws.Range("I55").Value = 231.5
$ws.Range("J55").Value = 1250.7646
$ws.Range("K55").Value = 231.5
$ws.Range("L55").Value = 1250.7646
$ws.Range("M55").Value = -17.5
$ws.Range("N55").Value = -1678.7646

# Row 64
$ws.Range("H64").Value = 5250
$ws.Range("I64").Value = 5000
$ws.Range("K64").Value = 5000
$ws.Range("M64").Value = -4752

# Row 67
$ws.Range("H67").Value = 5250
$ws.Range("I67").Value = 5000
$ws.Range("K67").Value = 5000
$ws.Range("M67").Value = -4142

# Row 76
$ws.Range("H76").Value = 5528
$ws.Range("J76").Value = 4500
$ws.Range("L76").Value = 4500
$ws.Range("N76").Value = -5130

# Row 79
$ws.Range("H79").Value = 5528
$ws.Range("J79").Value = 4500
$ws.Range("L79").Value = 4500
$ws.Range("N79").Value = -6684

# Row 98
$ws.Range("H98").Value = 978
$ws.Range("I98").Value = 978
$ws.Range("K98").Value = 978
$ws.Range("M98").Value = 520

# Row 122
$ws.Range("H122").Value = 978
$ws.Range("I122").Value = 978
$ws.Range("K122").Value = 2934
$ws.Range("M122").Value = -484

# Row 129
$ws.Range("H129").Value = 3164.5334
$ws.Range("I129").Value = 2872.125
$ws.Range("K129").Value = 8616.375
$ws.Range("M129").Value = -3616.375

# Row 134
$ws.Range("H134").Value = 50620
$ws.Range("J134").Value = 50620
$ws.Range("L134").Value = 50620
$ws.Range("N134").Value = -60760

# Row 138
$ws.Range("H138").Value = 3886
$ws.Range("I138").Value = 4249
$ws.Range("K138").Value = 12747
$ws.Range("M138").Value = -7607

# Row 141
$ws.Range("H141").Value = 3574.5
$ws.Range("I141").Value = 3574.5
$ws.Range("K141").Value = 10723.5
$ws.Range("M141").Value = -5543.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4561.15
$ws.Range("I32").Value = 4561.15
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4561.15
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4274.15
$ws.Range("N32").ClearContents()

# Row 61
$ws.Range("H61").Value = 3059.8
$ws.Range("I61").Value = 3059.8
$ws.Range("K61").Value = 3059.8
$ws.Range("M61").Value = -2847.8

# Row 122
$ws.Range("H122").Value = 2981.25
$ws.Range("I122").Value = 2184.75
$ws.Range("K122").Value = 6554.25
$ws.Range("M122").Value = -4104.25

# Row 136
$ws.Range("H136").Value = 3059.8
$ws.Range("I136").Value = 3059.8
$ws.Range("K136").Value = 9179.400000000001
$ws.Range("M136").Value = -6629.400000000001

# Row 23
$ws.Range("H23").Value = 3000
$ws.Range("J23").Value = 3000
$ws.Range("L23").Value = 3000
$ws.Range("N23").Value = -3566

$ws = $wb.Worksheets.Item("BSM")
# Row 31
$ws.Range("H31").Value = 5000
$ws.Range("J31").Value = 5000
$ws.Range("L31").Value = 5000
$ws.Range("N31").Value = -5504

# Row 99
$ws.Range("H99").Value = 3111
$ws.Range("I99").Value = 2222
$ws.Range("K99").Value = 2222
$ws.Range("M99").Value = -724

# Row 134
$ws.Range("H134").Value = 1293.4546
$ws.Range("I134").Value = 1382.9
$ws.Range("J134").Value = 399
$ws.Range("K134").Value = 4148.700000000001
$ws.Range("L134").Value = 1197
$ws.Range("M134").Value = -1613.700000000001
$ws.Range("N134").Value = -6267

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 4949.5
$ws.Range("I16").Value = 2000
$ws.Range("K16").Value = 2000
$ws.Range("M16").Value = -1713

# Row 53
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

# Row 58
$ws.Range("H58").Value = 1568.1786
$ws.Range("I58").Value = 1533.0385
$ws.Range("K58").Value = 1533.0385
$ws.Range("M58").Value = -1330.0385

# Row 62
$ws.Range("H62").Value = 2933.3333
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376

# Row 65
$ws.Range("H65").Value = 2933.3333
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880

# Row 86
$ws.Range("H86").Value = 9562.111000000001
$ws.Range("I86").Value = 10138.5
$ws.Range("K86").Value = 10138.5
$ws.Range("M86").Value = -9015.5

# Row 89
$ws.Range("H89").Value = 9562.111000000001
$ws.Range("I89").Value = 10138.5
$ws.Range("K89").Value = 50692.5
$ws.Range("M89").Value = -45076.5

# Row 105
$ws.Range("H105").Value = 2889.4285
$ws.Range("I105").Value = 2422.3333
$ws.Range("K105").Value = 2422.3333
$ws.Range("M105").Value = -675.3332999999998

# Row 113
$ws.Range("H113").Value = 4949.5
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 170

# Row 136
$ws.Range("H136").Value = 1568.1786
$ws.Range("I136").Value = 1533.0385
$ws.Range("K136").Value = 4599.1155
$ws.Range("M136").Value = -2049.1155

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 18333396
$ws.Range("I4").Value = 18333396
$ws.Range("K4").Value = 55000188
$ws.Range("M4").Value = -55000076

# Row 37
$ws.Range("H37").Value = 69992
$ws.Range("J37").Value = 69992
$ws.Range("L37").Value = 209976
$ws.Range("N37").Value = -210200

# Row 50
$ws.Range("H50").Value = 1068.6666
$ws.Range("J50").Value = 2499
$ws.Range("L50").Value = 7497
$ws.Range("N50").Value = -8459

# Row 53
$ws.Range("H53").Value = 1068.6666
$ws.Range("J53").Value = 2499
$ws.Range("L53").Value = 7497
$ws.Range("N53").Value = -8459

# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

# Row 140
$ws.Range("H140").Value = 10511.889
$ws.Range("J140").Value = 13815.077
$ws.Range("L140").Value = 41445.231
$ws.Range("N140").Value = -51805.231

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 8994
$ws.Range("I70").Value = 8499.5
$ws.Range("K70").Value = 8499.5
$ws.Range("M70").Value = -8229.5

# Row 73
$ws.Range("H73").Value = 8994
$ws.Range("I73").Value = 8499.5
$ws.Range("K73").Value = 8499.5
$ws.Range("M73").Value = -7563.5

# Row 122
$ws.Range("H122").Value = 3957.647
$ws.Range("I122").Value = 3759.4666
$ws.Range("J122").Value = 5444
$ws.Range("K122").Value = 11278.3998
$ws.Range("L122").Value = 16332
$ws.Range("M122").Value = -8828.399800000001
$ws.Range("N122").Value = -21232

# Row 26
$ws.Range("H26").Value = 3169.6667
$ws.Range("I26").Value = 2254.5
$ws.Range("J26").Value = 5000
$ws.Range("K26").Value = 2254.5
$ws.Range("L26").Value = 5000
$ws.Range("M26").Value = -1959.5
$ws.Range("N26").Value = -5590

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 6666
$ws.Range("I40").Value = 4999.5
$ws.Range("J40").Value = 9999
$ws.Range("K40").Value = 4999.5
$ws.Range("L40").Value = 9999
$ws.Range("M40").Value = -4863.5
$ws.Range("N40").Value = -10271

# Row 46
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()

# Row 132
$ws.Range("H132").Value = 2849.8333
$ws.Range("I132").Value = 2212.25
$ws.Range("J132").Value = 4125
$ws.Range("K132").Value = 6636.75
$ws.Range("L132").Value = 12375
$ws.Range("M132").Value = -4106.75
$ws.Range("N132").Value = -17435

# Row 136
$ws.Range("H136").Value = 2500.6667
$ws.Range("I136").Value = 1614.5714
$ws.Range("K136").Value = 4843.7142
$ws.Range("M136").Value = -2293.7142

# Row 140
$ws.Range("H140").Value = 75000
$ws.Range("J140").Value = 75000
$ws.Range("L140").Value = 75000
$ws.Range("N140").Value = -85360

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1667.0952
$ws.Range("I132").Value = 1704.5264
$ws.Range("K132").Value = 5113.5792
$ws.Range("M132").Value = -2583.5792

# Row 136
$ws.Range("H136").Value = 4842.7144
$ws.Range("I136").Value = 4983.1665
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 14949.4995
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -12399.4995
$ws.Range("N136").Value = -17100
